# Weekly update: insert a new week's price data for Vega Monumental
# Concepción - Acelga. Two new records (row "Primera" and "Segunda")
# are inserted right before the existing row 251, pushing every
# subsequent row down by two (the sheet grows from A1:R316 to A1:R318).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 251; this shifts rows
# 251-316 down to 253-318 and extends the sheet dimension accordingly.
$ws.Rows("251:252").Insert()

# New row 251: "Primera" quality record for the new date.
$ws.Cells.Item(251, 1).Value = 11
$ws.Cells.Item(251, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(251, 3).Value = "Bíobío"
$ws.Cells.Item(251, 4).Value = 44855
$ws.Cells.Item(251, 5).Value = 8
$ws.Cells.Item(251, 6).Value = 100112009
$ws.Cells.Item(251, 7).Value = "Acelga"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 200
$ws.Cells.Item(251, 11).Value = 700
$ws.Cells.Item(251, 12).Value = 800
$ws.Cells.Item(251, 13).Value = 750
$ws.Cells.Item(251, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(251, 15).Value = "Región de Ñuble"
$ws.Cells.Item(251, 16).Value = 750
$ws.Cells.Item(251, 17).Value = 1
$ws.Cells.Item(251, 18).Value = "Hortaliza"

# New row 252: "Segunda" quality record for the new date.
$ws.Cells.Item(252, 1).Value = 11
$ws.Cells.Item(252, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(252, 3).Value = "Bíobío"
$ws.Cells.Item(252, 4).Value = 44855
$ws.Cells.Item(252, 5).Value = 8
$ws.Cells.Item(252, 6).Value = 100112009
$ws.Cells.Item(252, 7).Value = "Acelga"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Segunda"
$ws.Cells.Item(252, 10).Value = 100
$ws.Cells.Item(252, 11).Value = 600
$ws.Cells.Item(252, 12).Value = 600
$ws.Cells.Item(252, 13).Value = 600
$ws.Cells.Item(252, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(252, 15).Value = "Región de Ñuble"
$ws.Cells.Item(252, 16).Value = 600
$ws.Cells.Item(252, 17).Value = 1
$ws.Cells.Item(252, 18).Value = "Hortaliza"
